$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ServiceImporting")

# Update the first (and only surviving) service row to "TestService3" and point its
# webservice path at the JIRA search JQL endpoint instead of issues.json.
$ws.Range("G2").Value = "/rest/api/2/search?jql=updated%20%3E%20startOfWeek(-2)%20ORDER%20BY%20updated%20DESC&startAt=0&maxResults=1500&fields=assignee,description,summary,created,updated,duedate,priority,status,worklog,key,id,project,timeestimate,timeoriginalestimate"
$ws.Range("B2").Value = "TestService3"

# Row now needs much more vertical room to show the long path text.
$ws.Rows.Item(2).RowHeight = 105

# Column G (Webservice Path) needs to be much wider to fit the long URL.
$ws.Columns.Item(7).ColumnWidth = 42.7109375

# The second test-service row (row 3, "TestService2") is no longer needed; deleting it
# shifts the trailing "W01" row up from row 4 to row 3.
$ws.Rows.Item(3).Delete()

# Rebuild the hyperlink list so only the remaining K2 cell (admin e-mail) is linked,
# then restore K2's original "Hyperlink" cell style (Add() re-applies it from scratch).
$ws.Range("K2").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("K2"), "mailto:123456@a")
$ws.Range("K2").Font.Underline = 2

# Move the active selection to D5, matching the saved view state.
$ws.Range("D5").Select()
